$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Arduino uno" entry from the materials list (row 10, column D)
$ws.Range("D10").ClearContents()

# Move the active selection to the now-empty cell, matching the recorded state
$ws.Range("D10").Select() | Out-Null
